$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained two new price records. They become the new rows 624 and
# 625; every existing record from row 624 downward (through the former last
# row 658) shifts down by two rows, ending at row 660. Inserting two whole
# rows at 624:625 reproduces that shift automatically, so only the two new
# rows need their values written explicitly.
$ws.Range("A624:R625").EntireRow.Insert()

# New row 624: Pepino ensalada, Primera, origin Provincia de Limarí
$ws.Range("A624").Value = 6
$ws.Range("B624").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C624").Value = 'Metropolitana'
$ws.Range("D624").Value = 45147
$ws.Range("E624").Value = 13
$ws.Range("F624").Value = 100112043
$ws.Range("G624").Value = 'Pepino ensalada'
$ws.Range("H624").Value = 'Sin especificar'
$ws.Range("I624").Value = 'Primera'
$ws.Range("J624").Value = 480
$ws.Range("K624").Value = 9000
$ws.Range("L624").Value = 10000
$ws.Range("M624").Value = 9479
$ws.Range("N624").Value = '$/caja 60 unidades'
$ws.Range("O624").Value = 'Provincia de Limarí'
$ws.Range("P624").Value = 158
$ws.Range("Q624").Value = 60
$ws.Range("R624").Value = 'Hortaliza'

# New row 625: Pepino ensalada, Segunda, origin Región de Arica y Parinacota
$ws.Range("A625").Value = 6
$ws.Range("B625").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C625").Value = 'Metropolitana'
$ws.Range("D625").Value = 45147
$ws.Range("E625").Value = 13
$ws.Range("F625").Value = 100112043
$ws.Range("G625").Value = 'Pepino ensalada'
$ws.Range("H625").Value = 'Sin especificar'
$ws.Range("I625").Value = 'Segunda'
$ws.Range("J625").Value = 120
$ws.Range("K625").Value = 7000
$ws.Range("L625").Value = 7000
$ws.Range("M625").Value = 7000
$ws.Range("N625").Value = '$/caja 80 unidades'
$ws.Range("O625").Value = 'Región de Arica y Parinacota'
$ws.Range("P625").Value = 88
$ws.Range("Q625").Value = 80
$ws.Range("R625").Value = 'Hortaliza'
